$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (weekly refresh: new Fecha + updated prices pulled in) ---

# Row 2
$ws.Range("D2").Value = 44405
$ws.Range("M2").Value = 50

# Row 4
$ws.Range("D4").Value = 44760
$ws.Range("N4:P4").Value = 2300
$ws.Range("S4").Value = 2300

# Row 5
$ws.Range("D5").Value = 44830
$ws.Range("M5").Value = 50
$ws.Range("N5:P5").Value = 2500
$ws.Range("S5").Value = 2500

# Row 6
$ws.Range("D6").Value = 44749
$ws.Range("M6").Value = 120
$ws.Range("N6:P6").Value = 2300
$ws.Range("S6").Value = 2300

# Row 7
$ws.Range("D7").Value = 44763
$ws.Range("M7").Value = 50
$ws.Range("N7:P7").Value = 2300
$ws.Range("S7").Value = 2300

# Row 8
$ws.Range("D8").Value = 44431
$ws.Range("M8").Value = 100
$ws.Range("N8:P8").Value = 1300
$ws.Range("S8").Value = 1300

# Row 9
$ws.Range("D9").Value = 45044
$ws.Range("M9").Value = 150
$ws.Range("N9:P9").Value = 3500
$ws.Range("S9").Value = 3500

# Row 10
$ws.Range("D10").Value = 44812
$ws.Range("M10").Value = 50
$ws.Range("N10:P10").Value = 2500
$ws.Range("S10").Value = 2500

# Row 11
$ws.Range("D11").Value = 45054
$ws.Range("M11").Value = 25

# Row 12
$ws.Range("D12").Value = 44476
$ws.Range("M12").Value = 80

# Row 13
$ws.Range("D13").Value = 44435
$ws.Range("M13").Value = 130

# Row 14
$ws.Range("D14").Value = 44357
$ws.Range("M14").Value = 35
$ws.Range("N14:P14").Value = 1000
$ws.Range("S14").Value = 1000

# Row 15
$ws.Range("D15").Value = 44432
$ws.Range("M15").Value = 30

# Row 16
$ws.Range("D16").Value = 44418
$ws.Range("M16").Value = 40
$ws.Range("N16:P16").Value = 1200
$ws.Range("S16").Value = 1200

# Row 17
$ws.Range("D17").Value = 45055
$ws.Range("M17").Value = 25
$ws.Range("N17:P17").Value = 2800
$ws.Range("S17").Value = 2800

# Row 18
$ws.Range("D18").Value = 44424
$ws.Range("M18").Value = 50

# Row 19
$ws.Range("D19").Value = 44438
$ws.Range("N19:P19").Value = 1200
$ws.Range("S19").Value = 1200

# Row 20
$ws.Range("D20").Value = 45041
$ws.Range("M20").Value = 80
$ws.Range("N20:P20").Value = 3500
$ws.Range("S20").Value = 3500

# Row 21
$ws.Range("D21").Value = 44343
$ws.Range("M21").Value = 60
$ws.Range("N21:P21").Value = 1300
$ws.Range("S21").Value = 1300

# Row 22
$ws.Range("D22").Value = 44811
$ws.Range("M22").Value = 60
$ws.Range("N22:P22").Value = 2500
$ws.Range("S22").Value = 2500

# Row 23
$ws.Range("D23").Value = 44753
$ws.Range("M23").Value = 160
$ws.Range("N23:P23").Value = 2300
$ws.Range("S23").Value = 2300

# Row 24
$ws.Range("D24").Value = 44762
$ws.Range("N24:P24").Value = 2300
$ws.Range("S24").Value = 2300

# Row 25
$ws.Range("D25").Value = 45042
$ws.Range("M25").Value = 25
$ws.Range("N25:P25").Value = 3500
$ws.Range("S25").Value = 3500

# Row 26
$ws.Range("D26").Value = 44417
$ws.Range("M26").Value = 80
$ws.Range("N26:P26").Value = 1200
$ws.Range("S26").Value = 1200

# --- Append new row 27 (same market/product constants as the rest of the sheet) ---
$ws.Range("A27").Value = 10
$ws.Range("B27").Value = "Vega Modelo de Temuco"
$ws.Range("C27").Value = "La Araucanía"
$ws.Range("D27").Value = 44473
$ws.Range("D27").NumberFormat = $ws.Range("D26").NumberFormat
$ws.Range("E27").Value = 9
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100108
$ws.Range("H27").Value = "Tropicales y subtropicales"
$ws.Range("I27").Value = 100108001
$ws.Range("J27").Value = "Guayaba"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 120
$ws.Range("N27:P27").Value = 1200
$ws.Range("Q27").Value = "`$/kilo"
$ws.Range("R27").Value = "Región de Arica y Parinacota"
$ws.Range("S27").Value = 1200
$ws.Range("T27").Value = 1
